$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "AIBT Courses Fees 2021.pdf"
$ws.Range("B2").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/onshore/coe/aibt/AIBT_Courses_Fees_2021_VOL_2.2.pdf"

# Update row 3 values
$ws.Range("A3").Value = "AIBT Onshore Q4 Promotion Brochure.pdf"
$ws.Range("B3").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/onshore/coe/aibt/AIBTOnshoreQ4Brochure_1OCT-31DEC21_VOL1.0.pdf"

# Remove rows 4-6 which are no longer present
$ws.Range("A4:B6").ClearContents()

# Match the author's final selection/scroll state
$ws.Range("A4:XFD6").Select()
